$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column M: a single-space value in every populated row (header through row 9) ---
$ws.Range("M1").Value = " "
$ws.Range("M2").Value = " "
$ws.Range("M3").Value = " "
$ws.Range("M4").Value = " "
$ws.Range("M5").Value = " "
$ws.Range("M6").Value = " "
$ws.Range("M7").Value = " "
$ws.Range("M8").Value = " "
$ws.Range("M9").Value = " "

# --- Row 3 (set first so its Column-B font change seeds the style; the rest of column B
#     copies that format instead of re-touching Font.Name, which keeps a single new cellXf) ---
$ws.Range("A3").Value = "MED-02"
$ws.Range("B3").Font.Name = "Aptos"
$ws.Range("B3").Value = "Usabilidad"
$ws.Range("C3").Value = "Protección frente a errores"
$ws.Range("D3").Value = "El error debe mostrar el campo equivocado"
$ws.Range("E3").Value = "Mensaje de error explicado (Si o No)"
$ws.Range("F3").Value = "Si"
$ws.Range("G3").Value = "Caso narrado"
$ws.Range("H3").Value = "Inspeccion"
$ws.Range("I3").Value = "No (Historia indica mensaje error)"
$ws.Range("J3").Value = "Baja"
$ws.Range("K3").Value = "P3"
$ws.Range("L3").Value = "El mensaje de error “401” no ofrece orientación al usuario sobre cómo proceder. Se sugiere incluir un texto explicativo o un enlace de ayuda para mejorar la experiencia de uso."

# --- Propagate the Column-B font (Aptos) down to B4:B9 via a format-only paste ---
$ws.Range("B3").Copy()
$ws.Range("B4:B9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Remaining data rows 4-9 ---
# Row 4
$ws.Range("A4").Value = "MED-03"
$ws.Range("B4").Value = "Seguridad"
$ws.Range("C4").Value = "Confidencialidad"
$ws.Range("D4").Value = "El cache debe ser borrado tras accion de login o sign up"
$ws.Range("E4").Value = "Replicaciòn de error (Se repite o no)"
$ws.Range("F4").Value = "No"
$ws.Range("G4").Value = "Caso narrado"
$ws.Range("H4").Value = "Prueba/Medición"
$ws.Range("I4").Value = "No (Historia indica almacenaminto)"
$ws.Range("J4").Value = "Bloqueante"
$ws.Range("K4").Value = "P1"
$ws.Range("L4").Value = "Se detectó que la sesión permanece activa tras cerrar la pestaña, lo que representa un riesgo de confidencialidad. Es necesario implementar cierre automático o invalidación del token."

# Row 5
$ws.Range("A5").Value = "MED-04"
$ws.Range("B5").Value = "Fiabilidad"
$ws.Range("C5").Value = "Recuperabilidad"
$ws.Range("D5").Value = "El turno debe ser recuperable fuera del search engine"
$ws.Range("E5").Value = "Almacenar datos en servidor (Refrescar pagina - Se recupera o no)"
$ws.Range("F5").Value = "Se recupera "
$ws.Range("G5").Value = "Caso narrado"
$ws.Range("H5").Value = "Prueba/Medición"
$ws.Range("I5").Value = "No (Histora indica perdida de datos)"
$ws.Range("J5").Value = "Media"
$ws.Range("K5").Value = "P2"
$ws.Range("L5").Value = "En caso de actualización del navegador, el sistema pierde la información del turno sin opción de recuperación. Se recomienda implementar una función de guardado temporal o recuperación."

# Row 6
$ws.Range("A6").Value = "MED-05"
$ws.Range("B6").Value = "Compatibilidad"
$ws.Range("C6").Value = "Interoperabilidad"
$ws.Range("D6").Value = "Tras exportar ambos calendarios deben ser identicos"
$ws.Range("E6").Value = "Comparar calendarios( Similitud >= 90)"
$ws.Range("F6").Value = ">= 90% en 3/3 mediciones"
$ws.Range("G6").Value = "Caso narrado"
$ws.Range("H6").Value = "Prueba/Medición"
$ws.Range("I6").Value = "No (Hisotria indica error de importacion)"
$ws.Range("J6").Value = "Media"
$ws.Range("K6").Value = "P3"
$ws.Range("L6").Value = "Al exportar el calendario, los datos no coinciden completamente con la aplicación original. Se debe revisar la interoperabilidad y el formato de exportación.."

# Row 7
$ws.Range("A7").Value = "MED-06"
$ws.Range("B7").Value = "Portabilidad"
$ws.Range("C7").Value = "Diseño Responsive"
$ws.Range("D7").Value = "El formulario debe ser visible en <= 80% en todo los dispositivos"
$ws.Range("E7").Value = "Visibilidad de formulario en >=80%"
$ws.Range("F7").Value = ">= 80% en 3/3 mediciones"
$ws.Range("G7").Value = "Caso narrado"
$ws.Range("H7").Value = "Prueba/Medición"
$ws.Range("I7").Value = "No (Hisotria indica error de UI)"
$ws.Range("J7").Value = "Media"
$ws.Range("K7").Value = "P3"
$ws.Range("L7").Value = "El formulario no mantiene su diseño correctamente en algunos dispositivos móviles. Se sugiere ajustar los estilos CSS y probar en diferentes resoluciones para asegurar el diseño responsive."

# Row 8
$ws.Range("A8").Value = "MED-07"
$ws.Range("B8").Value = "Usabilidad"
$ws.Range("C8").Value = "Estética de la UI"
$ws.Range("D8").Value = "El espacio entre ambos botones debe ser <= 200 pixeles"
$ws.Range("E8").Value = "Espacio de pixles >= 200p"
$ws.Range("F8").Value = ">= 200p en 3/3 mediciones"
$ws.Range("G8").Value = "Caso narrado"
$ws.Range("H8").Value = "Prueba/Medición"
$ws.Range("I8").Value = "No (Hisotria indica error de UI)"
$ws.Range("J8").Value = "Media"
$ws.Range("K8").Value = "P3"
$ws.Range("L8").Value = "Se observó que algunos espacios entre botones y campos de texto superan los 200 px, afectando la estética de la interfaz. Se recomienda unificar márgenes y espaciados."

# Row 9
$ws.Range("A9").Value = "MED-08"
$ws.Range("B9").Value = "Usabilidad"
$ws.Range("C9").Value = "Corrección funcional"
$ws.Range("D9").Value = "No debe haber errores duplicados"
$ws.Range("E9").Value = "Sin errores duplicados"
$ws.Range("F9").Value = "Turno unico"
$ws.Range("G9").Value = "Caso narrado"
$ws.Range("H9").Value = "Prueba/Medición"
$ws.Range("I9").Value = "No (Hisotria indica error de logic)"
$ws.Range("J9").Value = "Alta"
$ws.Range("K9").Value = "P1"
$ws.Range("L9").Value = "En la prueba de login se detectaron errores repetidos al ingresar credenciales, lo que afecta la funcionalidad. Se requiere depurar la validación de usuario y manejo de errores."

# --- Final selection, as left by the author ---
$ws.Range("M14").Select()
